# Append the new weekly-status row (row 7) to 工作表1, matching the
# columns already established by the existing rows:
#   A = 填寫日期 (date, plain number)
#   B = 姓名 (name)
#   C = 工作項目 (work item)
#   D = 工作完成目標 (completion target)
#   E = 工作完成現況 (current status)
#   F = 下周預計完成工作 (next week's planned work)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 12.2
$ws.Range("B7").Value = "溫舜元"
$ws.Range("C7").Value = "使用word2vec套件將文字轉向量，並計算頻率權重"
$ws.Range("D7").Value = "code完成"
$ws.Range("E7").Value = "初步code完成"
$ws.Range("F7").Value = "利用訓練的詞庫計算CNBC某些字詞出現的頻率"
